# Actualización automática 2025-06-30 13:15:09
#
# This script applies an incremental data refresh to the three report
# sheets (VENTAS POR GRUPO, VENTA MENSUAL, CUMPLIMIENTO MENSUAL) for
# GUERRERO FAREZ FABIAN MAURICIO: new sales picked up for two clients
# (ORTEGA ROMAN KLEBER ERWIN and PEREZ ROSALES EDGAR RICARDO) and all
# of the dependent totals / counters / ratios are refreshed to match.

$wb  = $excel.ActiveWorkbook
$wsGrupo  = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumpl   = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# ---------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": newly recorded sales by product group
# ---------------------------------------------------------------
# ORTEGA ROMAN KLEBER ERWIN (row 35) now has PANELES PVC sales
$wsGrupo.Range("Q35").Value = 417.24

# PEREZ ROSALES EDGAR RICARDO (row 40) now has FREGADEROS DE COCINA
# and GRIFERIAS sales
$wsGrupo.Range("E40").Value = 64.81999999999999
$wsGrupo.Range("G40").Value = 40.74

# Row 55 holds "N de 53" counters of how many of the 53 clients have
# a non-zero value in each column; the three columns above just moved
# from zero to non-zero, so their counters increment by one.
$wsGrupo.Range("E55").Value = "7 de 53"
$wsGrupo.Range("G55").Value = "4 de 53"
$wsGrupo.Range("Q55").Value = "3 de 53"

# ---------------------------------------------------------------
# Sheet "VENTA MENSUAL": junio column (F) totals per client, plus
# the grand total row 55
# ---------------------------------------------------------------
$wsMensual.Range("F35").Value = 3205.75
$wsMensual.Range("F40").Value = 1467.33
$wsMensual.Range("F55").Value = 80864.32000000001

# ---------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL": VENTA / POR CUMPLIR / CUMPLIMIENTO
# for the affected product groups, plus the TOTAL row
# ---------------------------------------------------------------
# FREGADEROS DE COCINA (row 4)
$wsCumpl.Range("D4").Value = 1762.62
$wsCumpl.Range("E4").Value = -759.6199999999999
$wsCumpl.Range("F4").Value = 1.757347956131605

# GRIFERIAS (row 6)
$wsCumpl.Range("D6").Value = 168.56
$wsCumpl.Range("E6").Value = -61.74000000000001
$wsCumpl.Range("F6").Value = 1.577981651376147

# PANELES PVC (row 14)
$wsCumpl.Range("D14").Value = 1087.7
$wsCumpl.Range("E14").Value = -121.7
$wsCumpl.Range("F14").Value = 1.125983436853002

# TOTAL (row 19)
$wsCumpl.Range("D19").Value = 83279.25
$wsCumpl.Range("E19").Value = 11168.19064517915
$wsCumpl.Range("F19").Value = 0.8817523209852145
